$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(3,7,11,15,19,23,27)
$cols = @("A","B","C","D","E")

foreach ($r1 in $rowPairs) {
    $r2 = $r1 + 1
    foreach ($c in $cols) {
        $v1 = $ws.Range("$c$r1").Value()
        $v2 = $ws.Range("$c$r2").Value()
        $ws.Range("$c$r1").Value = $v2
        $ws.Range("$c$r2").Value = $v1
    }
}

$ws.Range("F1:G29").Delete()
